# Update the "想去人数" (want-to-go count) figures in column F for the
# rows that changed between the previous and newly generated data pull.
#
# Sheet "展览" (exhibitions)
$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 213
$wsExpo.Range("F6").Value = 9805
$wsExpo.Range("F10").Value = 3388
$wsExpo.Range("F13").Value = 39
$wsExpo.Range("F18").Value = 266
$wsExpo.Range("F19").Value = 1432

# Sheet "全部类型" (all types) mirrors the same events at different rows
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 213
$wsAll.Range("F7").Value = 9805
$wsAll.Range("F11").Value = 3388
$wsAll.Range("F14").Value = 39
$wsAll.Range("F19").Value = 266
$wsAll.Range("F20").Value = 1432
